$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "model_3_9_11"
$ws.Range("B2").Value = 0.3122757150630314
$ws.Range("C2").Value = -0.3878164121035341
$ws.Range("D2").Value = 0.3033645620344083
$ws.Range("E2").Value = -0.3402954066405866
$ws.Range("F2").Value = 0.3147510232686105
$ws.Range("G2").Value = 0.4082628436808154
$ws.Range("H2").Value = 0.8238677727720853
$ws.Range("I2").Value = 0.6786833827987729
$ws.Range("J2").Value = 0.2815705006597134
$ws.Range("K2").Value = 0.4801269417292431
$ws.Range("L2").Value = 0.4153080131210624
$ws.Range("M2").Value = 0.6389544926525014
$ws.Range("N2").Value = -0.1789559170348032
$ws.Range("O2").Value = 0.666156108222859
$ws.Range("P2").Value = 21.7916881745875
$ws.Range("Q2").Value = 33.98044642326951

$ws.Range("A3").Value = "model_3_9_10"
$ws.Range("B3").Value = 0.3122333838625014
$ws.Range("C3").Value = -0.3881881382061314
$ws.Range("D3").Value = 0.3030575955119815
$ws.Range("E3").Value = -0.3388649417079264
$ws.Range("F3").Value = 0.3147520632138725
$ws.Range("G3").Value = 0.4082879733100626
$ws.Range("H3").Value = 0.8240884454443186
$ws.Range("I3").Value = 0.6789824388997013
$ws.Range("J3").Value = 0.2812699872614955
$ws.Range("K3").Value = 0.4801262130805984
$ws.Range("L3").Value = 0.4152732042656039
$ws.Range("M3").Value = 0.638974156997028
$ws.Range("N3").Value = -0.1790284848071404
$ws.Range("O3").Value = 0.666176609719243
$ws.Range("P3").Value = 21.79156507322635
$ws.Range("Q3").Value = 33.98032332190836

$ws.Range("A4").Value = "model_3_9_7"
$ws.Range("B4").Value = 0.3119060044201915
$ws.Range("C4").Value = -0.3885013932147794
$ws.Range("D4").Value = 0.3008139728822139
$ws.Range("E4").Value = -0.3285288928744834
$ws.Range("F4").Value = 0.314741788243231
$ws.Range("G4").Value = 0.4084823198890734
$ws.Range("H4").Value = 0.8242744071493637
$ws.Range("I4").Value = 0.6811682441474812
$ws.Range("J4").Value = 0.2790985805473813
$ws.Range("K4").Value = 0.4801334123474313
$ws.Range("L4").Value = 0.4150232962185295
$ws.Range("M4").Value = 0.6391262159300566
$ws.Range("N4").Value = -0.1795897067082433
$ws.Range("O4").Value = 0.666335142115856
$ws.Range("P4").Value = 21.79061329236803
$ws.Range("Q4").Value = 33.97937154105004

$ws.Range("A5").Value = "model_3_9_9"
$ws.Range("B5").Value = 0.3121004021848642
$ws.Range("C5").Value = -0.3885276410645728
$ws.Range("D5").Value = 0.3021192992921536
$ws.Range("E5").Value = -0.3345146434507751
$ws.Range("F5").Value = 0.3147519190915429
$ws.Range("G5").Value = 0.4083669169784174
$ws.Range("H5").Value = 0.8242899890068492
$ws.Range("I5").Value = 0.6798965555492933
$ws.Range("J5").Value = 0.2803560725735721
$ws.Range("K5").Value = 0.4801263140614327
$ws.Range("L5").Value = 0.4151674731282091
$ws.Range("M5").Value = 0.6390359277680853
$ws.Range("N5").Value = -0.1792564533973755
$ws.Range("O5").Value = 0.6662410101999072
$ws.Range("P5").Value = 21.79117840478396
$ws.Range("Q5").Value = 33.97993665346596

$ws.Range("A6").Value = "model_3_9_8"
$ws.Range("B6").Value = 0.3120054156176315
$ws.Range("C6").Value = -0.3885513830598575
$ws.Range("D6").Value = 0.3014748770439812
$ws.Range("E6").Value = -0.3315584456273644
$ws.Range("F6").Value = 0.3147470844170706
$ws.Range("G6").Value = 0.4084233051079328
$ws.Range("H6").Value = 0.8243040832808509
$ws.Range("I6").Value = 0.6805243712582065
$ws.Range("J6").Value = 0.2797350317962458
$ws.Range("K6").Value = 0.4801297015272261
$ws.Range("L6").Value = 0.4150968096045963
$ws.Range("M6").Value = 0.6390800459315975
$ws.Range("N6").Value = -0.1794192875126317
$ws.Range("O6").Value = 0.6662870065649774
$ws.Range("P6").Value = 21.79090225980533
$ws.Range("Q6").Value = 33.97966050848733

$ws.Range("A7").Value = "model_3_9_1"
$ws.Range("B7").Value = 0.3084882592109636
$ws.Range("C7").Value = -0.4118115955340795
$ws.Range("D7").Value = 0.283551577855195
$ws.Range("E7").Value = -0.2534033259490112
$ws.Range("F7").Value = 0.3140031278999948
$ws.Range("G7").Value = 0.4105112410841767
$ws.Range("H7").Value = 0.8381123501944095
$ws.Range("I7").Value = 0.6979857932034977
$ws.Range("J7").Value = 0.2633161318523062
$ws.Range("K7").Value = 0.480650962527902
$ws.Range("L7").Value = 0.4131047386095261
$ws.Range("M7").Value = 0.6407115115901826
$ws.Range("N7").Value = -0.1854486984954911
$ws.Range("O7").Value = 0.6679879270942479
$ws.Range("P7").Value = 21.78070393316217
$ws.Range("Q7").Value = 33.96946218184418

$ws.Range("A8").Value = "model_3_9_3"
$ws.Range("B8").Value = 0.3087347987058431
$ws.Range("C8").Value = -0.4120504818351249
$ws.Range("D8").Value = 0.2845790943533723
$ws.Range("E8").Value = -0.2576447375102697
$ws.Range("F8").Value = 0.3140816238837386
$ws.Range("G8").Value = 0.4103648845900645
$ws.Range("H8").Value = 0.838254163422061
$ws.Range("I8").Value = 0.6969847554513831
$ws.Range("J8").Value = 0.2642071715222852
$ws.Range("K8").Value = 0.4805959634868341
$ws.Range("L8").Value = 0.4132160151694273
$ws.Range("M8").Value = 0.6405972873733268
$ws.Range("N8").Value = -0.1850260593614117
$ws.Range("O8").Value = 0.6678688401160038
$ws.Range("P8").Value = 21.78141710530148
$ws.Range("Q8").Value = 33.97017535398349

$ws.Range("A9").Value = "model_3_9_0"
$ws.Range("B9").Value = 0.3083567290891852
$ws.Range("C9").Value = -0.4120814271405076
$ws.Range("D9").Value = 0.2830138728391041
$ws.Range("E9").Value = -0.2512070890539202
$ws.Range("F9").Value = 0.3139585549185038
$ws.Range("G9").Value = 0.4105893230462697
$ws.Range("H9").Value = 0.8382725338921037
$ws.Range("I9").Value = 0.6985096417466241
$ws.Range("J9").Value = 0.2628547443708195
$ws.Range("K9").Value = 0.4806821930587217
$ws.Range("L9").Value = 0.4130481296682532
$ws.Range("M9").Value = 0.6407724424834995
$ws.Range("N9").Value = -0.185674178704254
$ws.Range("O9").Value = 0.6680514519418377
$ws.Range("P9").Value = 21.7803235560647
$ws.Range("Q9").Value = 33.96908180474671

$ws.Range("A10").Value = "model_3_9_2"
$ws.Range("B10").Value = 0.308610564660307
$ws.Range("C10").Value = -0.4122531307310728
$ws.Range("D10").Value = 0.2840594464543116
$ws.Range("E10").Value = -0.2555015624200905
$ws.Range("F10").Value = 0.314041649789676
$ws.Range("G10").Value = 0.4104386352861261
$ws.Range("H10").Value = 0.8383744645606721
$ws.Range("I10").Value = 0.6974910122031623
$ws.Range("J10").Value = 0.2637569313139301
$ws.Range("K10").Value = 0.4806239717585462
$ws.Range("L10").Value = 0.4131603303103272
$ws.Range("M10").Value = 0.6406548487962346
$ws.Range("N10").Value = -0.1852390320109023
$ws.Range("O10").Value = 0.6679288520478533
$ws.Range("P10").Value = 21.78105769799157
$ws.Range("Q10").Value = 33.96981594667358

$ws.Range("A11").Value = "model_3_9_6"
$ws.Range("B11").Value = 0.3116146999171328
$ws.Range("C11").Value = -0.4155651058628855
$ws.Range("D11").Value = 0.2988598252808033
$ws.Range("E11").Value = -0.3192483615203974
$ws.Range("F11").Value = 0.3147745235895272
$ws.Range("G11").Value = 0.4086552508257884
$ws.Range("H11").Value = 0.8403405960687923
$ws.Range("I11").Value = 0.6830720340386278
$ws.Range("J11").Value = 0.2771489179231485
$ws.Range("K11").Value = 0.4801104759808881
$ws.Range("L11").Value = 0.4147619899893011
$ws.Range("M11").Value = 0.6392614886146266
$ws.Range("N11").Value = -0.1800890858563438
$ws.Range("O11").Value = 0.6664761736386612
$ws.Range("P11").Value = 21.78976677180408
$ws.Range("Q11").Value = 33.97852502048608

$ws.Range("A12").Value = "model_3_9_4"
$ws.Range("B12").Value = 0.3095389507590007
$ws.Range("C12").Value = -0.4196723660201345
$ws.Range("D12").Value = 0.2880785960966556
$ws.Range("E12").Value = -0.2722018704952491
$ws.Range("F12").Value = 0.3143322051345409
$ws.Range("G12").Value = 0.4098875051937502
$ws.Range("H12").Value = 0.8427788431225357
$ws.Range("I12").Value = 0.6935754346620506
$ws.Range("J12").Value = 0.267265347505313
$ws.Range("K12").Value = 0.4804203910836818
$ws.Range("L12").Value = 0.4135860821271153
$ws.Range("M12").Value = 0.6402245740314489
$ws.Range("N12").Value = -0.18364751298457
$ws.Range("O12").Value = 0.6674802595955392
$ws.Range("P12").Value = 21.78374506901683
$ws.Range("Q12").Value = 33.97250331769884

$ws.Range("A13").Value = "model_3_9_5"
$ws.Range("B13").Value = 0.3103102883591243
$ws.Range("C13").Value = -0.419838937363469
$ws.Range("D13").Value = 0.2917284440207519
$ws.Range("E13").Value = -0.2877694307442238
$ws.Range("F13").Value = 0.3145358309218028
$ws.Range("G13").Value = 0.4094296058742675
$ws.Range("H13").Value = 0.8428777270674466
$ws.Range("I13").Value = 0.6900196420611737
$ws.Range("J13").Value = 0.2705357949840079
$ws.Range("K13").Value = 0.4802777185225908
$ws.Range("L13").Value = 0.4139804190399427
$ws.Range("M13").Value = 0.6398668657418256
$ws.Range("N13").Value = -0.1823252199557868
$ws.Range("O13").Value = 0.6671073229234684
$ws.Range("P13").Value = 21.78598058618581
$ws.Range("Q13").Value = 33.97473883486781

$ws.Range("A14").Value = "model_3_9_12"
$ws.Range("B14").Value = 0.312452699094883
$ws.Range("C14").Value = -0.4257898026563882
$ws.Range("D14").Value = 0.3088504524693084
$ws.Range("E14").Value = -0.3761003462166852
$ws.Range("F14").Value = 0.3131971939069659
$ws.Range("G14").Value = 0.4081577783141961
$ws.Range("H14").Value = 0.8464104177693189
$ws.Range("I14").Value = 0.6733388618698716
$ws.Range("J14").Value = 0.2890924355351019
$ws.Range("K14").Value = 0.4812156487024868
$ws.Range("L14").Value = 0.4169460635904689
$ws.Range("M14").Value = 0.638872270735079
$ws.Range("N14").Value = -0.1786525158373433
$ws.Range("O14").Value = 0.6660703859481892
$ws.Range("P14").Value = 21.79220293555443
$ws.Range("Q14").Value = 33.98096118423643

$ws.Range("A15").Value = "model_3_9_13"
$ws.Range("B15").Value = 0.3125110213678282
$ws.Range("C15").Value = -0.4259047133279659
$ws.Range("D15").Value = 0.3093403339682986
$ws.Range("E15").Value = -0.3782859718514218
$ws.Range("F15").Value = 0.3132101096925528
$ws.Range("G15").Value = 0.408123155693585
$ws.Range("H15").Value = 0.8464786337078499
$ws.Range("I15").Value = 0.6728616044482909
$ws.Range("J15").Value = 0.2895515937931832
$ws.Range("K15").Value = 0.481206599120737
$ws.Range("L15").Value = 0.4169814805380874
$ws.Range("M15").Value = 0.6388451734916567
$ws.Range("N15").Value = -0.1785525347980088
$ws.Range("O15").Value = 0.6660421351190153
$ws.Range("P15").Value = 21.79237259587205
$ws.Range("Q15").Value = 33.98113084455405

$ws.Range("A16").Value = "model_3_9_20"
$ws.Range("B16").Value = 0.3134093715296185
$ws.Range("C16").Value = -0.4262981369567624
$ws.Range("D16").Value = 0.3197690845744839
$ws.Range("E16").Value = -0.4265120715632298
$ws.Range("F16").Value = 0.3132305344135111
$ws.Range("G16").Value = 0.4075898562308399
$ws.Range("H16").Value = 0.8467121869689194
$ws.Range("I16").Value = 0.6627015991513451
$ws.Range("J16").Value = 0.2996829774966867
$ws.Range("K16").Value = 0.4811922883240158
$ws.Range("L16").Value = 0.4177482059402328
$ws.Range("M16").Value = 0.6384276436925643
$ws.Range("N16").Value = -0.1770125059492254
$ws.Range("O16").Value = 0.6656068302119698
$ws.Range("P16").Value = 21.79498772904233
$ws.Range("Q16").Value = 33.98374597772433

$ws.Range("A17").Value = "model_3_9_21"
$ws.Range("B17").Value = 0.313416752105755
$ws.Range("C17").Value = -0.4263697377349864
$ws.Range("D17").Value = 0.3198989095564261
$ws.Range("E17").Value = -0.427124747059338
$ws.Range("F17").Value = 0.3132289416316554
$ws.Range("G17").Value = 0.4075854748020209
$ws.Range("H17").Value = 0.8467546922838658
$ws.Range("I17").Value = 0.6625751197144482
$ws.Range("J17").Value = 0.2998116889324839
$ws.Range("K17").Value = 0.481193404323466
$ws.Range("L17").Value = 0.4177570374952644
$ws.Range("M17").Value = 0.6384242122617382
$ws.Range("N17").Value = -0.1769998535329913
$ws.Range("O17").Value = 0.6656032526980292
$ws.Range("P17").Value = 21.7950092283623
$ws.Range("Q17").Value = 33.9837674770443

$ws.Range("A18").Value = "model_3_9_19"
$ws.Range("B18").Value = 0.313372348192944
$ws.Range("C18").Value = -0.4266558443434192
$ws.Range("D18").Value = 0.3191092038995366
$ws.Range("E18").Value = -0.4233662527007243
$ws.Range("F18").Value = 0.313243379979666
$ws.Range("G18").Value = 0.4076118348828147
$ws.Range("H18").Value = 0.8469245375257934
$ws.Range("I18").Value = 0.6633444749287024
$ws.Range("J18").Value = 0.2990221009558049
$ws.Range("K18").Value = 0.4811832879422537
$ws.Range("L18").Value = 0.4176994462850639
$ws.Range("M18").Value = 0.6384448565716656
$ws.Range("N18").Value = -0.1770759745263817
$ws.Range("O18").Value = 0.66562477587897
$ws.Range("P18").Value = 21.79487988504632
$ws.Range("Q18").Value = 33.98363813372833

$ws.Range("A19").Value = "model_3_9_15"
$ws.Range("B19").Value = 0.3131760985633187
$ws.Range("C19").Value = -0.4275666173562747
$ws.Range("D19").Value = 0.3161421262863434
$ws.Range("E19").Value = -0.4093070010286175
$ws.Range("F19").Value = 0.3132883115572969
$ws.Range("G19").Value = 0.4077283371405029
$ws.Range("H19").Value = 0.84746521172957
$ws.Range("I19").Value = 0.6662350919742968
$ws.Range("J19").Value = 0.2960685203402164
$ws.Range("K19").Value = 0.4811518061572565
$ws.Range("L19").Value = 0.4174731587519704
$ws.Range("M19").Value = 0.6385360891449307
$ws.Range("N19").Value = -0.177412402462882
$ws.Range("O19").Value = 0.6657198924117563
$ws.Range("P19").Value = 21.79430833338002
$ws.Range("Q19").Value = 33.98306658206202

$ws.Range("A20").Value = "model_3_9_22"
$ws.Range("B20").Value = 0.3134298640029689
$ws.Range("C20").Value = -0.4276015663341932
$ws.Range("D20").Value = 0.3201535233074218
$ws.Range("E20").Value = -0.4283471603740368
$ws.Range("F20").Value = 0.3132226951840009
$ws.Range("G20").Value = 0.4075776910134303
$ws.Range("H20").Value = 0.8474859589525797
$ws.Range("I20").Value = 0.6623270672720725
$ws.Range("J20").Value = 0.300068494654064
$ws.Range("K20").Value = 0.4811977809630683
$ws.Range("L20").Value = 0.4177762640585999
$ws.Range("M20").Value = 0.638418116138186
$ws.Range("N20").Value = -0.1769773759949105
$ws.Range("O20").Value = 0.6655968970498768
$ws.Range("P20").Value = 21.79504742335892
$ws.Range("Q20").Value = 33.98380567204093

$ws.Range("A21").Value = "model_3_9_16"
$ws.Range("B21").Value = 0.3131880251744185
$ws.Range("C21").Value = -0.4277638582197048
$ws.Range("D21").Value = 0.3163059053773124
$ws.Range("E21").Value = -0.4100825045187253
$ws.Range("F21").Value = 0.3132859138689543
$ws.Range("G21").Value = 0.4077212569889512
$ws.Range("H21").Value = 0.8475823024264637
$ws.Range("I21").Value = 0.6660755334140617
$ws.Range("J21").Value = 0.2962314388318348
$ws.Range("K21").Value = 0.4811534861229482
$ws.Range("L21").Value = 0.4174862523088159
$ws.Range("M21").Value = 0.6385305450712215
$ws.Range("N21").Value = -0.177391956843854
$ws.Range("O21").Value = 0.6657141123153516
$ws.Range("P21").Value = 21.79434306343128
$ws.Range("Q21").Value = 33.98310131211328

$ws.Range("A22").Value = "model_3_9_14"
$ws.Range("B22").Value = 0.3128866173315992
$ws.Range("C22").Value = -0.4283855835163604
$ws.Range("D22").Value = 0.3128273947723923
$ws.Range("E22").Value = -0.3940211371553084
$ws.Range("F22").Value = 0.3132754354353509
$ws.Range("G22").Value = 0.4079001857045897
$ws.Range("H22").Value = 0.8479513854196928
$ws.Range("I22").Value = 0.6694644039994326
$ws.Range("J22").Value = 0.2928572518970812
$ws.Range("K22").Value = 0.4811608279482569
$ws.Range("L22").Value = 0.4172326480535696
$ws.Range("M22").Value = 0.6386706394571381
$ws.Range("N22").Value = -0.1779086560029728
$ws.Range("O22").Value = 0.665860170809312
$ws.Range("P22").Value = 21.79346555477767
$ws.Range("Q22").Value = 33.98222380345968

$ws.Range("A23").Value = "model_3_9_23"
$ws.Range("B23").Value = 0.3134800022067934
$ws.Range("C23").Value = -0.4293962246051743
$ws.Range("D23").Value = 0.3212105048773288
$ws.Range("E23").Value = -0.4334400314552178
$ws.Range("F23").Value = 0.3131940282680461
$ws.Range("G23").Value = 0.4075479268097824
$ws.Range("H23").Value = 0.8485513456274347
$ws.Range("I23").Value = 0.6612973237529725
$ws.Range("J23").Value = 0.3011384097287697
$ws.Range("K23").Value = 0.4812178667408711
$ws.Range("L23").Value = 0.4178558454624702
$ws.Range("M23").Value = 0.6383948048110842
$ws.Range("N23").Value = -0.1768914247883542
$ws.Range("O23").Value = 0.6655725933113192
$ws.Range("P23").Value = 21.7951934828276
$ws.Range("Q23").Value = 33.98395173150961

$ws.Range("A24").Value = "model_3_9_18"
$ws.Range("B24").Value = 0.3133155615456715
$ws.Range("C24").Value = -0.429874748842703
$ws.Range("D24").Value = 0.3181267297479711
$ws.Range("E24").Value = -0.4186431118181593
$ws.Range("F24").Value = 0.3132684175528433
$ws.Range("G24").Value = 0.407645545889691
$ws.Range("H24").Value = 0.8488354182860025
$ws.Range("I24").Value = 0.6643016310599553
$ws.Range("J24").Value = 0.2980298591437379
$ws.Range("K24").Value = 0.4811657451018466
$ws.Range("L24").Value = 0.4176206241041449
$ws.Range("M24").Value = 0.6384712569017427
$ws.Range("N24").Value = -0.1771733230645631
$ws.Range("O24").Value = 0.6656523001257546
$ws.Range("P24").Value = 21.794714484486
$ws.Range("Q24").Value = 33.98347273316801

$ws.Range("A25").Value = "model_3_9_17"
$ws.Range("B25").Value = 0.313304851419365
$ws.Range("C25").Value = -0.430027859126817
$ws.Range("D25").Value = 0.3179746758346906
$ws.Range("E25").Value = -0.4179402298712926
$ws.Range("F25").Value = 0.3132680796961024
$ws.Range("G25").Value = 0.40765190388332
$ws.Range("H25").Value = 0.8489263111646724
$ws.Range("I25").Value = 0.6644497665962903
$ws.Range("J25").Value = 0.2978821970531992
$ws.Range("K25").Value = 0.4811659818247447
$ws.Range("L25").Value = 0.4176109495248493
$ws.Range("M25").Value = 0.6384762359581756
$ws.Range("N25").Value = -0.1771916832810885
$ws.Range("O25").Value = 0.665657491150928
$ws.Range("P25").Value = 21.79468329099396
$ws.Range("Q25").Value = 33.98344153967597

$ws.Range("A26").Value = "model_3_9_24"
$ws.Range("B26").Value = 0.3135272373363146
$ws.Range("C26").Value = -0.4310058384464885
$ws.Range("D26").Value = 0.3223810803045019
$ws.Range("E26").Value = -0.4391175997183465
$ws.Range("F26").Value = 0.3131566787293121
$ws.Range("G26").Value = 0.4075198859964473
$ws.Range("H26").Value = 0.8495068819353363
$ws.Range("I26").Value = 0.6601569136511622
$ws.Range("J26").Value = 0.302331158529115
$ws.Range("K26").Value = 0.4812440360901387
$ws.Range("L26").Value = 0.4179438081102746
$ws.Range("M26").Value = 0.6383728424646895
$ws.Range("N26").Value = -0.1768104502806034
$ws.Range("O26").Value = 0.6655496959823701
$ws.Range("P26").Value = 21.79533109500124
$ws.Range("Q26").Value = 33.98408934368324
